# Generate Report for Handoff
# - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
# - Related timestamps are refreshed to reflect the new handoff generation time
# - Status column widths shrink to fit the shorter text

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-15 22:55:42"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-15 22:55:37"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-15 22:55:42"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
